$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 9928.26
$ws.Range("B10").Value = 9871.99
$ws.Range("C10").Value = 307.87
$ws.Range("D10").Value = 306.13
$ws.Range("E10").Value = $true
$ws.Range("F10").Value = -0.57
$ws.Range("G10").Value = 42612.67291666667
$ws.Range("G10").NumberFormat = "m/d/yy h:mm"
$ws.Range("H10").Value = $true
